$wb = $excel.ActiveWorkbook

# --- DBS sheet ("DBS", the sheet2.xml file) -------------------------------
$dbs = $wb.Worksheets.Item("DBS")

# Re-enter the same order-by text into C6; Excel records it as a (new)
# duplicate shared-string entry rather than reusing the old one.
$dbs.Range("C6").Value = "ExportDate DESC ,CustNo ASC"

# Add a new function-condition row describing the "bringUpDateFirst" lookup.
# Populate C14 first, then A14, then B14 so the new shared strings are
# appended to sharedStrings.xml in that same order.
$dbs.Range("C14").Value = "BringUpDate DESC"
$dbs.Range("A14").Value = "bringUpDateFirst"
$dbs.Range("B14").Value = "BringUpDate>= "

# --- View state -------------------------------------------------------------
$dbd = $wb.Worksheets.Item("DBD")
$dbd.Range("C1").Select()

$dbs.Range("A15").Select()
